# Quantum Health Vendor Template - "Implementing edits to template page"
#
# Functional change in the target diff:
#   - Sheet "QH-Vendor_form_template", cell A1 header text changes from
#     "Name" to "Vendor" (the other header cells B1:E1 - Description, url,
#     phoneNumber, SSO - are unchanged; they only shift shared-string
#     indices as a natural side effect of "Name" being removed from the
#     shared strings table and "Vendor" being appended).
#   - The worksheet's saved cursor/selection moves from B6 to B17.
#
# (The Descriptions sheet's cell text is untouched by this edit - its
# shared-string indices shift too, but purely as a side effect of the
# A1 edit above, and that happens automatically once the shared string
# table is updated.)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("QH-Vendor_form_template")
$ws2 = $wb.Worksheets.Item("Descriptions")

# --- Header cell edit: "Name" -> "Vendor" ---
$ws1.Range("A1").Value = "Vendor"

# --- Restore/update the saved selection on the template sheet ---
$ws1.Activate() | Out-Null
$ws1.Range("B17").Select() | Out-Null

# Descriptions sheet keeps its existing selection (D27); nothing to change there.
